# Evaluation Form update: switch curriculum from Electrical Engineering to
# Civil Engineering, fix the student name typo, and replace the subject
# list (rows 6-11) with the Civil Engineering subjects. Rows 12-15 (the
# old extra GE/NSTP/PATHFIT rows) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block -----------------------------------------------------
$ws.Range("E2").Value = "Civil Engineering"
$ws.Range("B2").Value = "Test Student1"
$ws.Range("B3").Value = "18-0000"

# --- Subject rows (6-11) become the Civil Engineering curriculum -----
$ws.Range("A6").Value = 63
$ws.Range("B6").Value = "CE 24"
$ws.Range("C6").Value = "CE Projec 2"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "CE 18"
$ws.Range("F6").Value = 1.25
$ws.Range("G6").Value = "Passed"

$ws.Range("A7").Value = 64
$ws.Range("B7").Value = "CE 25"
$ws.Range("C7").Value = "Elective 3: Computer Software in Structural Analysis"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "CE 11 / CE 12"
$ws.Range("F7").Value = 1.5
$ws.Range("G7").Value = "Passed"

$ws.Range("A8").Value = 65
$ws.Range("B8").Value = "CE 26"
$ws.Range("C8").Value = "Elective 4: Prestresssed Concrete Design"
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "CE 12"
$ws.Range("F8").Value = 1.75
$ws.Range("G8").Value = "Passed"

$ws.Range("A9").Value = 66
$ws.Range("B9").Value = "CE 27"
$ws.Range("C9").Value = "CE Elective 5"
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = " "
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = "Passed"

$ws.Range("A10").Value = 67
$ws.Range("B10").Value = "CE 28"
$ws.Range("C10").Value = "CE Elective 6"
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = " "
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "Passed"

$ws.Range("A11").Value = 68
$ws.Range("B11").Value = "CE 29"
$ws.Range("C11").Value = "CE Integration Course 2"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "CE 23"
$ws.Range("F11").Value = 1.5
$ws.Range("G11").Value = "Passed"

# --- Drop the now-unused rows 12-15 -----------------------------------
$ws.Rows("12:15").Delete()

# --- Keep the selection in sync with the shrunk table -----------------
[void]$ws.Range("G6:G11").Select()
